# Update "想去人数" (number of people interested) counts on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, as published by
# the latest scrape (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value  = 164
$wsExhibit.Range("F4").Value  = 124
$wsExhibit.Range("F5").Value  = 1270
$wsExhibit.Range("F6").Value  = 17794
$wsExhibit.Range("F8").Value  = 247
$wsExhibit.Range("F10").Value = 6720
$wsExhibit.Range("F12").Value = 152
$wsExhibit.Range("F17").Value = 144
$wsExhibit.Range("F19").Value = 184
$wsExhibit.Range("F27").Value = 106
$wsExhibit.Range("F28").Value = 5138
$wsExhibit.Range("F31").Value = 11913
$wsExhibit.Range("F36").Value = 3906

# --- Sheet 4: 全部类型 ---------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value  = 164
$wsAll.Range("F4").Value  = 124
$wsAll.Range("F5").Value  = 1270
$wsAll.Range("F6").Value  = 17794
$wsAll.Range("F8").Value  = 247
$wsAll.Range("F10").Value = 6720
$wsAll.Range("F12").Value = 152
$wsAll.Range("F17").Value = 144
$wsAll.Range("F19").Value = 184
$wsAll.Range("F27").Value = 106
$wsAll.Range("F28").Value = 5138
$wsAll.Range("F33").Value = 11913
$wsAll.Range("F38").Value = 3906
